$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns F, G, H (row 1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header formatting (bold, centered, bordered) from the existing
# header cell E1 onto the new header cells so the style matches the other
# headers exactly (reuses the same style index instead of creating a new one).
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the header text after the paste-special (paste-special with
# formats-only does not touch values, but keep this explicit/safe).
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New boolean (FALSE) data cells for rows 2-5 in columns F, G, H
$ws.Range("F2:H5").Value = $false
